$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generic Cluster RA")

# Add explanatory notes in column G next to the assumption/parameter rows
# (written in the same order the author typed them, so new shared-string
# entries land in the same order as the target workbook)
$ws.Range("G8").Value = '<-- uses t distribution with DF correction'
$ws.Range("G7").Value = '<-- Should almost always be "2"'
$ws.Range("G6").Value = '<-- Should almost always be 0.05'
$ws.Range("G5").Value = '<-- Should almost always be 0.80'
$ws.Range("G9").Value = '<-- Make this specific to the study context'
$ws.Range("G10").Value = '<-- Make this specific to the study context'
$ws.Range("G11").Value = '<-- Your best estimate based on your study design'
$ws.Range("G12").Value = '<-- Your best estimate based on your study design'
$ws.Range("G13").Value = '<-- Your best estimate based on your study design'

# Reflect the selection that was active when the author saved the file
$ws.Range("G14").Select()
